$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Bruselas (repollito)" at the top
# of the data (row 19, right after the existing row 18), pushing every
# subsequent record down by one row (old row 19 -> 20, ..., old row 80 -> 81).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record's data. The
# non-numeric / descriptive columns repeat the same values used by the
# surrounding rows for this market/category.
$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 44715
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 100112035
$ws.Cells.Item(19, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 12
$ws.Cells.Item(19, 11).Value = 28000
$ws.Cells.Item(19, 12).Value = 28000
$ws.Cells.Item(19, 13).Value = 28000
$ws.Cells.Item(19, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(19, 16).Value = 2800
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = "Hortaliza"
